# Apply updated odds values to sheet1 (Jogos_do_Dia_Betfair_Back_Lay_2025-12-25)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "I2" = 2.26
    "J2" = 3.2
    "L2" = 1.45
    "T2" = 1.84
    "V2" = 1.79
    "T3" = 1.85
    "P4" = 1.56
    "Q4" = 2
    "G5" = 1.7
    "I5" = 8
    "J5" = 3.75
    "P5" = 1.92
    "V5" = 1.15
    "W5" = 2.44
    "F6" = 3.15
    "G6" = 3.45
    "H6" = 2.32
    "I6" = 2.56
    "J6" = 3.35
    "P6" = 1.9
    "Q6" = 1.9
    "R6" = 1.36
    "V6" = 1.64
    "W6" = 1.41
    "Z6" = 17
    "AB6" = 15.5
    "AD6" = 12
    "AG6" = 15
    "AJ6" = 65
    "AN6" = 42
    "AO6" = 24
    "F7" = 1.43
    "G7" = 1.5
    "H7" = 8
    "I7" = 9.6
    "J7" = 4.5
    "K7" = 5.2
    "M7" = 1.05
    "P7" = 2.2
    "Q7" = 1.67
    "S7" = 2.68
    "T7" = 1.89
    "U7" = 1.94
    "V7" = 1.11
    "AC7" = 11.5
    "AM7" = 150
    "G8" = 1.69
    "H8" = 6.6
    "J8" = 3.25
    "K8" = 4.3
    "L8" = 1.45
    "N8" = 2.56
    "O8" = 1.49
    "P8" = 1.54
    "Q8" = 2.26
    "R8" = 1.19
    "S8" = 4.9
    "T8" = 2.36
    "U8" = 1.58
    "V8" = 1.11
    "W8" = 2.44
    "I9" = 5.2
    "J9" = 2.82
    "N9" = 2.74
    "O9" = 1.26
    "P9" = 1.98
    "Q9" = 1.66
    "R9" = 1.39
    "T9" = 1.48
    "U9" = 1.81
    "V9" = 1.27
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

